# Re-order three groups of duplicate Artfynd observation rows (3/4, 9/10,
# 18/19/20). The underlying records were re-sequenced by the data export;
# only the cells that actually differ between the old and new row contents
# are touched so that untouched cells (dates, shared columns, etc.) keep
# their original representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3 <-> 4 -----------------------------------------------------
$ws.Range("A3").Value2 = 131039759
$ws.Range("B3").Value2 = 91829
$ws.Range("E3").Value2 = 5442
$ws.Range("F3").Value2 = "Tallticka"
$ws.Range("G3").Value2 = "Porodaedalea pini"
$ws.Range("H3").Value2 = "(Brot.) Murrill"
$ws.Range("P3").Value2 = "Gotvad, Dlr"
$ws.Range("Q3").Value2 = 479059
$ws.Range("R3").Value2 = 6792254
$ws.Range("S3").Value2 = 10
$ws.Range("AC3").Value2 = ""

$ws.Range("A4").Value2 = 131041641
$ws.Range("B4").Value2 = 79243
$ws.Range("E4").Value2 = 6425
$ws.Range("F4").Value2 = "Garnlav"
$ws.Range("G4").Value2 = "Alectoria sarmentosa"
$ws.Range("H4").Value2 = "(Ach.) Ach."
$ws.Range("P4").Value2 = "Tandbergsvasseln, Dlr"
$ws.Range("Q4").Value2 = 479078
$ws.Range("R4").Value2 = 6791615
$ws.Range("S4").Value2 = 50
$ws.Range("AC4").Value2 = "Rikligt i en radie av ca 50 meter, synfältet"

# --- Rows 9 <-> 10 ------------------------------------------------------
$ws.Range("A9").Value2 = 131040374
$ws.Range("B9").Value2 = 79001
$ws.Range("E9").Value2 = 228912
$ws.Range("F9").Value2 = "Mörk kolflarnlav"
$ws.Range("G9").Value2 = "Carbonicola myrmecina"
$ws.Range("H9").Value2 = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q9").Value2 = 479088
$ws.Range("R9").Value2 = 6792211

$ws.Range("A10").Value2 = 131039523
$ws.Range("B10").Value2 = 79243
$ws.Range("E10").Value2 = 6425
$ws.Range("F10").Value2 = "Garnlav"
$ws.Range("G10").Value2 = "Alectoria sarmentosa"
$ws.Range("H10").Value2 = "(Ach.) Ach."
$ws.Range("Q10").Value2 = 479079
$ws.Range("R10").Value2 = 6792517

# --- Rows 18 -> 19 -> 20 -> 18 (cyclic rotation) ------------------------
$ws.Range("A18").Value2 = 131041965
$ws.Range("B18").Value2 = 57884
$ws.Range("E18").Value2 = 100109
$ws.Range("F18").Value2 = "Tretåig hackspett"
$ws.Range("G18").Value2 = "Picoides tridactylus"
$ws.Range("H18").Value2 = "(Linnaeus, 1758)"
$ws.Range("M18").Value2 = "färska spår"
$ws.Range("P18").Value2 = "Tandbergsvasseln, Dlr"
$ws.Range("Q18").Value2 = 479096
$ws.Range("R18").Value2 = 6792085
$ws.Range("S18").Value2 = 10
$ws.Range("AC18").Value2 = ""

$ws.Range("A19").Value2 = 131039828
$ws.Range("M19").Value2 = "bobygge"
$ws.Range("P19").Value2 = "Gotvad, Dlr"
$ws.Range("Q19").Value2 = 479059
$ws.Range("R19").Value2 = 6792254
$ws.Range("AE19").Value2 = $true

$ws.Range("A20").Value2 = 131039579
$ws.Range("B20").Value2 = 79243
$ws.Range("E20").Value2 = 6425
$ws.Range("F20").Value2 = "Garnlav"
$ws.Range("G20").Value2 = "Alectoria sarmentosa"
$ws.Range("H20").Value2 = "(Ach.) Ach."
$ws.Range("M20").Value2 = ""
$ws.Range("Q20").Value2 = 479079
$ws.Range("R20").Value2 = 6792475
$ws.Range("S20").Value2 = 50
$ws.Range("AC20").Value2 = "Rikligt till måttligt i en radie av ca 50 meter, synfältet"
$ws.Range("AE20").Value2 = $false
